$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="28.123.79"},
    @{Cell="E2"; Value="  -4.55%  "},
    @{Cell="D3"; Value="1.833.00"},
    @{Cell="E3"; Value="  -3.05%  "},
    @{Cell="D4"; Value="1.002"},
    @{Cell="E4"; Value="  -0.46%  "},
    @{Cell="D5"; Value="329.30"},
    @{Cell="E5"; Value="  -3.10%  "},
    @{Cell="E6"; Value="  -0.39%  "},
    @{Cell="D7"; Value="0.4656"},
    @{Cell="E7"; Value="  -1.94%  "},
    @{Cell="E8"; Value="  -3.05%  "},
    @{Cell="D9"; Value="46.29"},
    @{Cell="D10"; Value="0.07892"},
    @{Cell="E10"; Value="  -1.63%  "},
    @{Cell="D11"; Value="0.9612"},
    @{Cell="E11"; Value="  -2.90%  "},
    @{Cell="D12"; Value="22.00"},
    @{Cell="E12"; Value="  -4.52%  "},
    @{Cell="D13"; Value="1.807.11"},
    @{Cell="E13"; Value="  -3.27%  "},
    @{Cell="D14"; Value="5.674"},
    @{Cell="E14"; Value="  -4.40%  "},
    @{Cell="D15"; Value="6.908"},
    @{Cell="E15"; Value="  -2.52%  "},
    @{Cell="D16"; Value="0.06856"},
    @{Cell="E16"; Value="  +1.24%  "},
    @{Cell="E17"; Value="  -0.53%  "},
    @{Cell="D18"; Value="86.82"},
    @{Cell="E18"; Value="  -2.47%  "},
    @{Cell="D19"; Value="0.000010000"},
    @{Cell="E19"; Value="  -1.93%  "},
    @{Cell="D20"; Value="16.68"},
    @{Cell="E20"; Value="  -3.56%  "},
    @{Cell="E21"; Value="  -0.29%  "},
    @{Cell="D22"; Value="28.144.56"},
    @{Cell="E22"; Value="  -4.51%  "},
    @{Cell="D23"; Value="5.338"},
    @{Cell="E23"; Value="  -2.90%  "},
    @{Cell="E24"; Value="  -5.43%  "},
    @{Cell="D25"; Value="2.094"},
    @{Cell="E25"; Value="  -2.61%  "},
    @{Cell="D26"; Value="2.042.65"},
    @{Cell="E26"; Value="  -2.48%  "},
    @{Cell="D27"; Value="152.79"},
    @{Cell="E27"; Value="  -2.95%  "},
    @{Cell="D28"; Value="19.28"},
    @{Cell="E28"; Value="  -1.82%  "},
    @{Cell="D29"; Value="5.780"},
    @{Cell="E29"; Value="  -10.90%  "},
    @{Cell="D30"; Value="1.976"},
    @{Cell="E30"; Value="  -3.42%  "},
    @{Cell="D31"; Value="117.40"},
    @{Cell="E31"; Value="  -1.18%  "},
    @{Cell="B32"; Value="ImmutableX"},
    @{Cell="C32"; Value="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"},
    @{Cell="D32"; Value="0.9379"},
    @{Cell="E32"; Value="  -5.90%  "},
    @{Cell="B33"; Value="Stellar"},
    @{Cell="C33"; Value="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"},
    @{Cell="D33"; Value="0.09275"},
    @{Cell="E33"; Value="  -2.60%  "},
    @{Cell="D34"; Value="5.308"},
    @{Cell="E34"; Value="  -3.05%  "},
    @{Cell="E35"; Value="  -4.33%  "},
    @{Cell="D36"; Value="3.359"},
    @{Cell="E36"; Value="  -5.03%  "},
    @{Cell="D37"; Value="0.05942"},
    @{Cell="E37"; Value="  -7.36%  "},
    @{Cell="D38"; Value="0.02153"},
    @{Cell="E38"; Value="  -3.88%  "},
    @{Cell="D39"; Value="1.150"},
    @{Cell="E39"; Value="  -4.07%  "},
    @{Cell="E40"; Value="  -0.37%  "},
    @{Cell="D41"; Value="0.5607"},
    @{Cell="E41"; Value="  -3.81%  "},
    @{Cell="D42"; Value="9.939"},
    @{Cell="E42"; Value="  -5.46%  "},
    @{Cell="D43"; Value="0.1776"},
    @{Cell="E43"; Value="  -2.46%  "},
    @{Cell="E44"; Value="  -2.70%  "},
    @{Cell="D45"; Value="2.213"},
    @{Cell="E45"; Value="  -8.26%  "},
    @{Cell="D46"; Value="11.69"},
    @{Cell="E46"; Value="  -3.55%  "},
    @{Cell="D47"; Value="0.5283"},
    @{Cell="E47"; Value="  -3.92%  "},
    @{Cell="D48"; Value="0.07053"},
    @{Cell="E48"; Value="  -3.83%  "},
    @{Cell="D49"; Value="1.840"},
    @{Cell="E49"; Value="  -5.66%  "},
    @{Cell="D50"; Value="112.11"},
    @{Cell="E50"; Value="  -3.49%  "},
    @{Cell="D51"; Value="1.000"},
    @{Cell="E51"; Value="  -0.64%  "}
)

foreach ($chg in $changes) {
    $cell = $ws.Range($chg.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
    $cell.Style = "Normal"
}
